# Clear the sample "Name"/"Email" row (A2:B2), including the mailto: hyperlink
# on B2, while preserving B2's existing (Hyperlink) cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove any hyperlinks attached to the sheet (the mailto: link on B2)
$ws.Hyperlinks.Delete()

# Clear the cell contents of A2 and B2 (drops the now-unused shared strings
# "Zemu" / "zemu@gmail.com"), but keep B2's cell formatting/style intact.
$ws.Range("A2").ClearContents()
$ws.Range("B2").ClearContents()

# Update the active selection to match: whole row 2 selected, A2 active.
$ws.Range("A2:XFD2").Select()
